$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora) changes from 7 to 8 for all data rows (2-51)
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "8"

# Column D (Price) and E (Volume(1h)) updates
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.38%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.92%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.016"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.14%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07523"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.87%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.598"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.97%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9174"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.10%"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.85%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1178"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.10%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1823"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.29%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08985"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.73%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04102"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.72%"
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.52%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001287"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.68%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005792"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.74%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.339"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.377"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.01%"
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.44%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.274"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.74%"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1352"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.84%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3224"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "16.23%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04095"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.93%"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.30%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003896"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "7.17%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.34%"
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02392"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.74%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05194"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.45%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006309"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.81%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007785"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.49%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1324"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.88%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007399"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.78%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006933"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.55%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3244"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.53%"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006587"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.71%"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04544"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.62%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004206"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.09%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
